$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the medal table into chronological order (oldest Olympiad first) ---
# Snapshot the current data rows (2-7, columns A-H) before they get overwritten.
$data = @()
for ($r = 2; $r -le 7; $r++) {
    $row = @()
    for ($c = 1; $c -le 8; $c++) {
        $row += $ws.Cells.Item($r, $c).Value2
    }
    $data += ,$row
}

# Write the rows back in reverse order.
for ($i = 0; $i -lt 6; $i++) {
    $srcRow = $data[5 - $i]
    $destR = 2 + $i
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($destR, $c).Value = $srcRow[$c - 1]
    }
}

# --- Formatting tweak: bold the column headers in row 1 ---
$ws.Range("B1:H1").Font.Bold = $true

# --- Move the active selection (was C9, now C10) ---
$null = $ws.Range("C10").Select()

# --- Page setup fix ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
